# Updated cryptos list on Fri Mar 17 03:52:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (the sheet stores prices like "25.689.75" /
# "1.000" / "0.9983" as plain strings, not numbers) without leaving a
# lingering custom cell style behind once we're done.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- Simple price / volume refreshes (coin identity unchanged) ---
Set-TextValue "D2" "25.689.75"
$ws.Range("E2").Value = "  +5.37%  "

Set-TextValue "D3" "1.703.00"
$ws.Range("E3").Value = "  +3.16%  "

Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.52%  "

Set-TextValue "D5" "330.63"
$ws.Range("E5").Value = "  +6.33%  "

Set-TextValue "D6" "0.9983"
$ws.Range("E6").Value = "  -0.34%  "

Set-TextValue "D7" "0.3685"
$ws.Range("E7").Value = "  +1.05%  "

$ws.Range("E8").Value = "  +3.17%  "

$ws.Range("E10").Value = "  +3.87%  "

Set-TextValue "D11" "0.07331"
$ws.Range("E11").Value = "  +4.42%  "

Set-TextValue "D12" "0.9990"
$ws.Range("E12").Value = "  -0.35%  "

Set-TextValue "D13" "6.186"
$ws.Range("E13").Value = "  +4.39%  "

Set-TextValue "D14" "19.90"
$ws.Range("E14").Value = "  +2.91%  "

Set-TextValue "D15" "6.822"
$ws.Range("E15").Value = "  +3.62%  "

Set-TextValue "D16" "1.698.42"
$ws.Range("E16").Value = "  +2.73%  "

Set-TextValue "D17" "0.00001069"
$ws.Range("E17").Value = "  +2.70%  "

Set-TextValue "D18" "0.06622"
$ws.Range("E18").Value = "  +0.40%  "

Set-TextValue "D19" "80.80"
$ws.Range("E19").Value = "  +3.68%  "

Set-TextValue "D20" "0.9973"
$ws.Range("E20").Value = "  -0.38%  "

Set-TextValue "D21" "16.11"
$ws.Range("E21").Value = "  +3.47%  "

Set-TextValue "D22" "6.026"
$ws.Range("E22").Value = "  +1.84%  "

$ws.Range("E23").Value = "  +4.36%  "

Set-TextValue "D24" "25.653.00"
$ws.Range("E24").Value = "  +5.26%  "

Set-TextValue "D25" "2.449"
$ws.Range("E25").Value = "  -1.25%  "

Set-TextValue "D26" "2.473"
$ws.Range("E26").Value = "  +6.04%  "

Set-TextValue "D27" "149.23"
$ws.Range("E27").Value = "  +1.45%  "

Set-TextValue "D28" "19.12"
$ws.Range("E28").Value = "  +3.13%  "

$ws.Range("E29").Value = "  +8.63%  "

Set-TextValue "D30" "1.890.62"
$ws.Range("E30").Value = "  +2.98%  "

Set-TextValue "D31" "127.87"
$ws.Range("E31").Value = "  +3.37%  "

Set-TextValue "D32" "4.091"
$ws.Range("E32").Value = "  +0.56%  "

Set-TextValue "D33" "5.938"
$ws.Range("E33").Value = "  +5.53%  "

Set-TextValue "D34" "0.08474"

Set-TextValue "D35" "1.684"
$ws.Range("E35").Value = "  +0.95%  "

Set-TextValue "D36" "12.68"
$ws.Range("E36").Value = "  +3.51%  "

Set-TextValue "D37" "5.299"
$ws.Range("E37").Value = "  +2.51%  "

Set-TextValue "D38" "1.269"
$ws.Range("E38").Value = "  +2.13%  "

Set-TextValue "D39" "0.06214"
$ws.Range("E39").Value = "  +3.28%  "

# --- Rows 40-45: three coin pairs swapped rank, so name/link/price/volume
#     all move together. Write the whole block (B:E) for rows 40-45. ---
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "8.490"
$ws.Range("E40").Value = "  +4.58%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D41" "0.2117"
$ws.Range("E41").Value = "  +2.75%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.02250"
$ws.Range("E42").Value = "  +1.81%  "

Set-TextValue "D43" "0.6089"
$ws.Range("E43").Value = "  +3.76%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "14.08"
$ws.Range("E44").Value = "  +12.09%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D45" "0.9978"
$ws.Range("E45").Value = "  -0.34%  "

# --- remaining simple refreshes ---
Set-TextValue "D46" "3.837"
$ws.Range("E46").Value = "  +1.80%  "

Set-TextValue "D47" "0.5830"
$ws.Range("E47").Value = "  +4.11%  "

Set-TextValue "D48" "125.64"
$ws.Range("E48").Value = "  +2.73%  "

Set-TextValue "D49" "1.999"
$ws.Range("E49").Value = "  +3.06%  "

Set-TextValue "D50" "0.07204"
$ws.Range("E50").Value = "  +4.63%  "

$ws.Range("E51").Value = "  +3.22%  "
